# 2024.12.18 LeetCodeHot100 三刷 栈Done
# Update the "last submitted" date, "times solved" count, difficulty rating,
# and notes for the stack-related problems that were reviewed for the third
# time (三刷) on 2024.12.18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 34: 155. 最小栈 (Min Stack)
$ws.Range("D34").Value = "2024.12.18"
$ws.Range("F34").Value = 3

# Row 46: 20. 有效的括号 (Valid Parentheses)
$ws.Range("D46").Value = "2024.12.18"
$ws.Range("F46").Value = 3

# Row 105: 394. 字符串解码 (Decode String) - also bump the star rating
$ws.Range("C105").Value = "⭐⭐⭐⭐⭐⭐"
$ws.Range("D105").Value = "2024.12.18"
$ws.Range("F105").Value = 3

# Row 140: 739. 每日温度 (Daily Temperatures) - update the note
$ws.Range("D140").Value = "2024.12.18"
$ws.Range("G140").Value = "while、temperatures[st.peek()]下标"

# Row 149: 84. 柱状图中最大的矩形 (Largest Rectangle in Histogram)
$ws.Range("D149").Value = "2024.12.18"
$ws.Range("F149").Value = 3

# Reflect where the user left the selection / viewport after the edits
$ws.Range("A99").Select()
$ws.Range("C105").Select()
